# Update countries & provincias Spain
# - Reorders three country rows (Cabo Verde before Zimbabue, Granada before
#   Santa Sede, Islas Malvinas before Montserrat) by swapping the full data
#   row (country name + stats) between the two row positions involved.
# - Refreshes the COVID-19 stat columns (Casos totales/Nuevos casos/Casos
#   activos/Recuperados/Muertes hoy/Muertes) for a number of countries with
#   the newer snapshot's figures.
# - Updates the "Datos actualizados..." timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row-order swaps: move the whole data row so the country list order
#    changes while each country keeps (and, where noted, updates) its own
#    statistics.
# ---------------------------------------------------------------------

# Rows 118/119: Cabo Verde now comes before Zimbabue.
$ws.Range("A118").Value = "Cabo Verde"
$ws.Range("B118").Value = 8322
$ws.Range("C118").Value = 124
$ws.Range("D118").Value = 7234
$ws.Range("E118").Value = 994
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 94

$ws.Range("A119").Value = "Zimbabue"
$ws.Range("B119").Value = 8269
$ws.Range("C119").Value = 12
$ws.Range("D119").Value = 7785
$ws.Range("E119").Value = 248
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 236

# Rows 208/209: Granada now comes before Santa Sede.
$ws.Range("A208").Value = "Granada"
$ws.Range("B208").Value = 28
$ws.Range("C208").Value = 1
$ws.Range("D208").Value = 24
$ws.Range("E208").Value = 4
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("B209").Value = 27
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 15
$ws.Range("E209").Value = 12
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# Rows 216/217: Islas Malvinas now comes before Montserrat.
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0

$ws.Range("A217").Value = "Montserrat"
$ws.Range("B217").Value = 13
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 12
$ws.Range("E217").Value = 0
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 1

# ---------------------------------------------------------------------
# 2) Plain statistic refreshes (no reordering involved).
# ---------------------------------------------------------------------

$ws.Range("B4").Value = 8823897
$ws.Range("C4").Value = 76614
$ws.Range("D4").Value = 5735723
$ws.Range("E4").Value = 2858129
$ws.Range("G4").Value = 761
$ws.Range("H4").Value = 230045

$ws.Range("B6").Value = 5381224
$ws.Range("C6").Value = 25574
$ws.Range("D6").Value = 4817898
$ws.Range("E6").Value = 406400
$ws.Range("G6").Value = 398
$ws.Range("H6").Value = 156926

$ws.Range("B10").Value = 1081336
$ws.Range("C10").Value = 11968
$ws.Range("D10").Value = 881113
$ws.Range("E10").Value = 171610
$ws.Range("G10").Value = 275
$ws.Range("H10").Value = 28613

$ws.Range("B11").Value = 1007711
$ws.Range("C11").Value = 8769
$ws.Range("D11").Value = 907379
$ws.Range("E11").Value = 70332
$ws.Range("G11").Value = 198
$ws.Range("H11").Value = 30000

$ws.Range("B12").Value = 886214
$ws.Range("C12").Value = 3098
$ws.Range("D12").Value = 803846
$ws.Range("E12").Value = 48273
$ws.Range("G12").Value = 62
$ws.Range("H12").Value = 34095

$ws.Range("B20").Value = 427808
$ws.Range("C20").Value = 10458
$ws.Range("E20").Value = 103597

$ws.Range("B31").Value = 250797
$ws.Range("C31").Value = 12474
$ws.Range("D31").Value = 94902
$ws.Range("E31").Value = 153818
$ws.Range("G31").Value = 106
$ws.Range("H31").Value = 2077

$ws.Range("B40").Value = 128515
$ws.Range("C40").Value = 649
$ws.Range("D40").Value = 104562
$ws.Range("E40").Value = 21325
$ws.Range("G40").Value = 6
$ws.Range("H40").Value = 2628

$ws.Range("B48").Value = 106397
$ws.Range("C48").Value = 167
$ws.Range("D48").Value = 98813
$ws.Range("E48").Value = 1397
$ws.Range("G48").Value = 11
$ws.Range("H48").Value = 6187

$ws.Range("B52").Value = 95835
$ws.Range("C52").Value = 697
$ws.Range("D52").Value = 88787
$ws.Range("E52").Value = 5342
$ws.Range("G52").Value = 12
$ws.Range("H52").Value = 1706

$ws.Range("B64").Value = 61930
$ws.Range("C64").Value = 48
$ws.Range("D64").Value = 57285
$ws.Range("E64").Value = 3516

$ws.Range("B65").Value = 59043
$ws.Range("C65").Value = 784
$ws.Range("D65").Value = 39214
$ws.Range("E65").Value = 18536
$ws.Range("G65").Value = 15
$ws.Range("H65").Value = 1293

$ws.Range("B97").Value = 17749
$ws.Range("C97").Value = 217
$ws.Range("E97").Value = 5607

$ws.Range("B98").Value = 16968
$ws.Range("C98").Value = 158
$ws.Range("D98").Value = 16301
$ws.Range("E98").Value = 423
$ws.Range("G98").Value = 6
$ws.Range("H98").Value = 244

$ws.Range("B155").Value = 2807
$ws.Range("C155").Value = 48
$ws.Range("D155").Value = 2301
$ws.Range("E155").Value = 453

$ws.Range("D199").Value = 65
$ws.Range("E199").Value = 8

$ws.Range("B204").Value = 38
$ws.Range("C204").Value = 1
$ws.Range("E204").Value = 9

# ---------------------------------------------------------------------
# 3) Timestamp update.
# ---------------------------------------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 25 de Octubre de 2020 a las 01:18"
